# Insert a new data row at row 122 (pushes existing rows 122..165 down to 123..166)
# and populate it with the new Locoto "Tercera" quality record dated 45009.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 122, shifting rows 122:165 down to 123:166.
$ws.Range("A122:R122").EntireRow.Insert()

# Fill the newly inserted row 122 with the new record's values.
$ws.Cells.Item(122, 1).Value  = 1
$ws.Cells.Item(122, 2).Value  = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(122, 3).Value  = "Arica y Parinacota"
$ws.Cells.Item(122, 4).Value  = 45009
$ws.Cells.Item(122, 5).Value  = 15
$ws.Cells.Item(122, 6).Value  = 100112042
$ws.Cells.Item(122, 7).Value  = "Locoto"
$ws.Cells.Item(122, 8).Value  = "Sin especificar"
$ws.Cells.Item(122, 9).Value  = "Tercera"
$ws.Cells.Item(122, 10).Value = 150
$ws.Cells.Item(122, 11).Value = 43000
$ws.Cells.Item(122, 12).Value = 45000
$ws.Cells.Item(122, 13).Value = 44333
$ws.Cells.Item(122, 14).Value = "$/caja 20 kilos"
$ws.Cells.Item(122, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(122, 16).Value = 2217
$ws.Cells.Item(122, 17).Value = 20
$ws.Cells.Item(122, 18).Value = "Hortaliza"
